# New crime data collected - weekly CompStat update (49th Precinct)
# Week-over-week rollover: Volume/Number and report date range bump,
# plus refreshed weekly/28-day/YTD crime statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (Volume/Number banner + reporting week range) ---
$ws.Range("A8").Value = "Volume 31   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/28/2024  Through  11/3/2024"

# --- Row 14 (Murder) ---
$ws.Range("N14").Value = -75

# --- Row 15 (Rape) ---
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 31
$ws.Range("K15").Value = 47.619047619047
$ws.Range("L15").Value = 106.666666666667
$ws.Range("M15").Value = 72.222222222222
$ws.Range("N15").Value = 6.896551724137

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = -71.428571428571
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -16
$ws.Range("I16").Value = 212
$ws.Range("J16").Value = 244
$ws.Range("K16").Value = -13.114754098360
$ws.Range("L16").Value = -5.357142857142
$ws.Range("M16").Value = -19.083969465648
$ws.Range("N16").Value = -63.636363636363

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = -27.027027027027
$ws.Range("I17").Value = 397
$ws.Range("J17").Value = 353
$ws.Range("K17").Value = 12.464589235127
$ws.Range("L17").Value = 28.064516129032
$ws.Range("M17").Value = 81.278538812785
$ws.Range("N17").Value = 38.327526132404

# --- Row 18 (Burglary) ---
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 11.111111111111
$ws.Range("I18").Value = 140
$ws.Range("J18").Value = 206
$ws.Range("K18").Value = -32.038834951456
$ws.Range("L18").Value = 23.893805309734
$ws.Range("M18").Value = -56.790123456790
$ws.Range("N18").Value = -88.682295877122

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -35.294117647058
$ws.Range("F19").Value = 61
$ws.Range("H19").Value = 5.172413793103
$ws.Range("I19").Value = 727
$ws.Range("J19").Value = 567
$ws.Range("K19").Value = 28.218694885361
$ws.Range("L19").Value = 39.272030651341
$ws.Range("M19").Value = 90.813648293963
$ws.Range("N19").Value = 43.960396039604

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -41.666666666666
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 42
$ws.Range("H20").Value = -19.047619047619
$ws.Range("I20").Value = 388
$ws.Range("J20").Value = 444
$ws.Range("K20").Value = -12.612612612612
$ws.Range("L20").Value = 29.765886287625
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -76.004947433518

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 46
$ws.Range("E21").Value = -34.782608695652
$ws.Range("F21").Value = 155
$ws.Range("G21").Value = 172
$ws.Range("H21").Value = -9.883720930232
$ws.Range("I21").Value = 1898
$ws.Range("J21").Value = 1840
$ws.Range("K21").Value = 3.152173913043
$ws.Range("L21").Value = 27.382550335570
$ws.Range("M21").Value = 34.801136363636
$ws.Range("N21").Value = -55.550351288056

# --- Row 22 (Transit) ---
$ws.Range("F22").Value = "0"
$ws.Range("H22").Value = -100

# --- Row 23 (Housing) ---
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -12.5
$ws.Range("I23").Value = 93
$ws.Range("J23").Value = 103
$ws.Range("K23").Value = -9.708737864077
$ws.Range("L23").Value = -1.063829787234
$ws.Range("M23").Value = 66.071428571428

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -45.833333333333
$ws.Range("F24").Value = 86
$ws.Range("G24").Value = 141
$ws.Range("H24").Value = -39.007092198581
$ws.Range("I24").Value = 1067
$ws.Range("J24").Value = 1373
$ws.Range("K24").Value = -22.286962855061
$ws.Range("L24").Value = -5.575221238938
$ws.Range("M24").Value = 34.892541087231

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -83.333333333333
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = -40.579710144927
$ws.Range("I25").Value = 408
$ws.Range("J25").Value = 592
$ws.Range("K25").Value = -31.081081081081
$ws.Range("L25").Value = -8.314606741573

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 37.5
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 36
$ws.Range("H26").Value = 11.111111111111
$ws.Range("I26").Value = 479
$ws.Range("J26").Value = 447
$ws.Range("K26").Value = 7.158836689038
$ws.Range("L26").Value = 5.739514348785
$ws.Range("M26").Value = -11.296296296296

# --- Row 27 (UCR Rape*) ---
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 35
$ws.Range("K27").Value = 12.903225806451
$ws.Range("L27").Value = 2.941176470588

# --- Row 28 (Other Sex Crimes) ---
$ws.Range("C28").Value = "0"
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 200
$ws.Range("L28").Value = 23.404255319148

# --- Row 31 (Hate Crimes) ---
$ws.Range("D31").Value = "0"
$ws.Range("E31").Value = "***.*"
